# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-03-28 Thursday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-03-29 Friday", 2)

# Update the answer table, cell by cell (row, column), to avoid ambiguity
# between duplicate original values that map to different new values.
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "16÷3=5, 1" },
    @{ Row = 1;  Col = 2; Text = "77÷3=25, 2" },
    @{ Row = 1;  Col = 3; Text = "32÷3=10, 2" },
    @{ Row = 1;  Col = 4; Text = "81÷3=27, 0" },
    @{ Row = 1;  Col = 5; Text = "51÷2=25, 1" },

    @{ Row = 5;  Col = 1; Text = "17÷7=2, 3" },
    @{ Row = 5;  Col = 2; Text = "69÷8=8, 5" },
    @{ Row = 5;  Col = 3; Text = "85÷5=17, 0" },
    @{ Row = 5;  Col = 4; Text = "95÷2=47, 1" },
    @{ Row = 5;  Col = 5; Text = "98÷3=32, 2" },

    @{ Row = 9;  Col = 1; Text = "11÷6=1, 5" },
    @{ Row = 9;  Col = 2; Text = "95÷2=47, 1" },
    @{ Row = 9;  Col = 3; Text = "72÷5=14, 2" },
    @{ Row = 9;  Col = 4; Text = "61÷4=15, 1" },
    @{ Row = 9;  Col = 5; Text = "98÷2=49, 0" },

    @{ Row = 13; Col = 1; Text = "22÷4=5, 2" },
    @{ Row = 13; Col = 2; Text = "12÷5=2, 2" },
    @{ Row = 13; Col = 3; Text = "16÷7=2, 2" },
    @{ Row = 13; Col = 4; Text = "46÷9=5, 1" },
    @{ Row = 13; Col = 5; Text = "59÷9=6, 5" },

    @{ Row = 17; Col = 1; Text = "39÷6=6, 3" },
    @{ Row = 17; Col = 2; Text = "34÷8=4, 2" },
    @{ Row = 17; Col = 3; Text = "19÷7=2, 5" },
    @{ Row = 17; Col = 4; Text = "64÷9=7, 1" },
    @{ Row = 17; Col = 5; Text = "69÷8=8, 5" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
